$wb = $excel.ActiveWorkbook

# The "Chart" sheet holds a daily GSC breadcrumbs export (Date | Invalid |
# Valid? columns, one row per day). This commit refreshes the export: the
# oldest date (2025-10-30, sheet row 2) has rolled off the window, so that
# row is deleted outright. Excel shifts every later row (and its "Items"
# count in column C) up by one, which reproduces the newly pulled data for
# 2025-10-31 .. 2026-01-26 exactly, drops the table down to 89 data+header
# rows, and removes the now-unused "2025-10-30" shared string (cascading
# the shared-string index shift into the other two sheets' header rows).
$chart = $wb.Worksheets.Item("Chart")
$chart.Rows.Item(2).Delete()
